$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "ParticipantsTab" query text (cell B2) was rewritten to a new Neo4j
# Cypher query (broader OPTIONAL MATCH based traversal, sorted sample ids,
# ORDER BY p.participant_id). Assigning the new text automatically causes
# the shared-string table to drop the old, now-unused string and append
# the new one, which naturally renumbers the other shared query-string
# references used by row 2/3/4 (B/C columns) exactly like the target diff.
$newQuery = @'
MATCH (p:participant)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
OPTIONAL MATCH (p)<--(diag:diagnosis)
OPTIONAL MATCH (samp)<--(f:file)
OPTIONAL MATCH (f)<--(g:genomic_info)
WITH s, p, samp, f, g, diag
WHERE g.platform in ['Illumina HiSeq X Ten']
with p
OPTIONAL MATCH (p)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
WITH s, p, apoc.coll.sort(collect(distinct samp.sample_id)) as samp
RETURN
coalesce(p.participant_id,'') as `Participant ID`,
coalesce(s.study_name, '') as `Study Name`,
coalesce(s.phs_accession,'') as `Accession`,
coalesce(p.gender,'') as `Gender`,
coalesce(apoc.text.join(samp, ','), '') as `Samples`
ORDER BY p.participant_id LIMIT 100
'@

$ws.Range("B2").Value = $newQuery

# The longer replacement text wraps across more lines, so the row grows
# taller (186 -> 279, i.e. 18 lines @ 15.5pt/line instead of 12).
$ws.Rows(2).RowHeight = 279

# Match the author's final cursor/selection position in the saved file.
$ws.Range("B5").Select() | Out-Null
